# Re-generate the statistics with fixed minutes and seconds formatting
# in the haul field (column D): zero-pad single-digit minutes and seconds,
# e.g. "17 ч. 1 мин. 54 сек." -> "17 ч. 01 мин. 54 сек."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2

    if ($val -ne $null -and $val -match '^(\d+) ч\. (\d+) мин\. (\d+) сек\.$') {
        $hours = $Matches[1]
        $minutes = $Matches[2].PadLeft(2, '0')
        $seconds = $Matches[3].PadLeft(2, '0')
        $newVal = "$hours ч. $minutes мин. $seconds сек."

        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}
